$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, applied while preserving text (string) cell type
$updates = [ordered]@{
    'E2' = '3.86%'
    'D3' = '28.44'
    'E3' = '4.75%'
    'D4' = '4.938'
    'E4' = '1.72%'
    'D5' = '0.06532'
    'E5' = '2.11%'
    'D6' = '7.250'
    'E6' = '4.33%'
    'D7' = '1.365'
    'E7' = '14.84%'
    'D8' = '0.9122'
    'E8' = '4.24%'
    'D9' = '0.1559'
    'E9' = '3.10%'
    'D10' = '0.06833'
    'E10' = '34.84%'
    'D11' = '0.07659'
    'E11' = '1.83%'
    'D12' = '0.02986'
    'E12' = '1.03%'
    'D13' = '0.08974'
    'E13' = '-0.03%'
    'D14' = '0.001596'
    'E14' = '2.17%'
    'D15' = '0.0006530'
    'E15' = '2.34%'
    'D16' = '0.006038'
    'E16' = '-2.24%'
    'E17' = '-0.23%'
    'D18' = '3.391'
    'E18' = '2.54%'
    'D19' = '2.244'
    'E19' = '-0.37%'
    'D20' = '0.3151'
    'E20' = '0.47%'
    'D21' = '0.1351'
    'E21' = '0.22%'
    'D22' = '4.028'
    'E22' = '3.15%'
    'D24' = '0.04478'
    'E24' = '1.50%'
    'D25' = '0.001188'
    'E25' = '1.24%'
    'E26' = '12.28%'
    'E28' = '-1.65%'
    'D29' = '0.0001638'
    'E29' = '-15.62%'
    'D40' = '0.04169'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '0.006697'
    'E41' = '-1.64%'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D42' = '0.1234'
    'E42' = '5.16%'
    'D43' = '0.002171'
    'E43' = '-0.86%'
    'D44' = '0.01242'
    'D45' = '0.00005558'
    'E45' = '7.04%'
    'D46' = '1.561'
    'E46' = '-7.48%'
    'D47' = '0.01852'
    'E47' = '0.16%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number format so Excel does not auto-convert the
    # numeric-looking / percent-looking string into a real number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default cell style so no stray formatting is introduced.
    $cell.Style = "Normal"
}
